$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.641.78"
$ws.Range("E2").Value = "  +0.77%  "
$ws.Range("D3").Value = "3.641.70"
$ws.Range("E3").Value = "  +5.91%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.24"
$ws.Range("E5").Value = "  -1.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.13"
$ws.Range("E6").Value = "  -1.10%  "
$ws.Range("D7").Value = "3.630.44"
$ws.Range("E7").Value = "  +5.85%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.612"
$ws.Range("E8").Value = "  +1.75%  "
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.197"
$ws.Range("E10").Value = "  -4.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.84"
$ws.Range("E11").Value = "  +22.77%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.603"
$ws.Range("E12").Value = "  +2.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "48.46"
$ws.Range("E13").Value = "  -0.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000288"
$ws.Range("E14").Value = "  +0.75%  "
$ws.Range("D15").Value = "4.228.12"
$ws.Range("E15").Value = "  +6.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "668.44"
$ws.Range("E16").Value = "  -3.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.87"
$ws.Range("E17").Value = "  +2.09%  "
$ws.Range("D18").Value = "3.651.87"
$ws.Range("E18").Value = "  +6.50%  "
$ws.Range("D19").Value = "70.743.79"
$ws.Range("E19").Value = "  +1.02%  "
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.75"
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.40"
$ws.Range("E22").Value = "  -0.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.928"
$ws.Range("E23").Value = "  +2.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.08"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "100.36"
$ws.Range("E25").Value = "  -0.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.90"
$ws.Range("E26").Value = "  -0.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.78"
$ws.Range("E27").Value = "  +3.99%  "
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.11%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.97"
$ws.Range("E29").Value = "  +3.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.80"
$ws.Range("E30").Value = "  +3.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.35"
$ws.Range("E31").Value = "  +0.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.01"
$ws.Range("E32").Value = "  +2.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.40"
$ws.Range("E33").Value = "  -3.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.28"
$ws.Range("E34").Value = "  +1.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.00"
$ws.Range("E35").Value = "  +4.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "581.81"
$ws.Range("E36").Value = "  +1.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.03"
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("E38").Value = "  +3.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "58.27"
$ws.Range("E39").Value = "  -0.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.25%  "
$ws.Range("D41").Value = "3.581.06"
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0454"
$ws.Range("E42").Value = "  +7.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.141"
$ws.Range("E43").Value = "  +1.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.343"
$ws.Range("E44").Value = "  +2.53%  "
$ws.Range("B45").Value = "PEPE"
$ws.Range("C45").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D45").Value = "0.0₃0743"
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "34.68"
$ws.Range("E46").Value = "  -2.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.70"
$ws.Range("E47").Value = "  +0.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.87"
$ws.Range("E48").Value = "  +6.70%  "
$ws.Range("E49").Value = "  +2.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.92"
$ws.Range("E50").Value = "  +1.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.95"
$ws.Range("E51").Value = "  +7.98%  "
